$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three input values in row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 3000

# Update the selected cell on the sheet
$ws.Range("D6").Select()
